$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "cryptoAAVEpolygon.xlsx"
$ws.Range("A2").Value = "cryptoAerodromeBase.xlsx"
$ws.Range("A3").Value = "cryptoChainlinkPolygon.xlsx"
$ws.Range("A4").Value = "cryptoDogeBnb.xlsx"
$ws.Range("A5").Value = "cryptoMorphoBase.xlsx"
$ws.Range("A6").Value = "cryptoMystPolygon.xlsx"
$ws.Range("A7").Value = "cryptoPaxgoldPolygon.xlsx"
$ws.Range("A8").Value = "cryptoSolanaPolygon.xlsx"
$ws.Range("A9").Value = "cryptoWrappedBTCPolygon.xlsx"

$ws.Range("A10:A18").ClearContents()
$ws.Range("A10:A40").Select()
